$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 128 (shifting the existing "view" rows down by one)
# so the data for the new "modify_location" stored procedure can be added.
$ws.Rows.Item(128).Insert()

$ws.Range("A128").Value = "stored procedure"
$ws.Range("B128").Value = "modify_location"
$ws.Range("C128").Value = "validate user's inputs, insert a new location or update an old one"

$ws.Range("D131").Select()
